# Apply the "Use RDString for some settings messages" change:
# 1. Replace hard-coded example names in a few translation strings with
#    {0}/{1} placeholders (DisableEffects, HideUiElements, JudgmentVisuals).
# 2. Add new translation rows for a hitsound-volume setting
#    (Miscellaneous sheet: SET_HITSOUND_VOLUME / CURRENT_HITSOUND_VOLUME).

$wb = $excel.ActiveWorkbook

# --- DisableEffects sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("DisableEffects")

# FILTER row
$ws.Cells.Item(4, 2).Value = "Disable VFX filters ({0}, {1}, etc.)"
$ws.Cells.Item(4, 3).Value = "VFX 필터 ({0}, {1}, 등..) 끄기"
$ws.Cells.Item(4, 4).Value = "Desactivar efectos de filtro ({0}, {1}, etc.)"

# HALL_OF_MIRRORS row
$ws.Cells.Item(7, 2).Value = 'Disable "{0}" effect'
$ws.Cells.Item(7, 3).Value = '"{0}" 이펙트 끄기'
$ws.Cells.Item(7, 4).Value = 'Desactivar "{0}"'

# --- HideUiElements sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("HideUiElements")

# JUDGE_TEXT row
$ws.Cells.Item(5, 2).Value = "Hide judgment text ({0}, {1}, etc.)"
$ws.Cells.Item(5, 3).Value = "판정 텍스트 숨기기 ({0}, {1}, 등..)"
$ws.Cells.Item(5, 4).Value = "Ocultar texto de juicios ({0}, {1}, etc.)"

# --- JudgmentVisuals sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("JudgmentVisuals")

# HIDE_PERFECTS row
$ws.Cells.Item(11, 2).Value = 'Hide "{0}" judgments'
$ws.Cells.Item(11, 3).Value = '"{0}" 판정 숨기기'
$ws.Cells.Item(11, 4).Value = 'Esconder juicios de "{0}"'

# --- Miscellaneous sheet: add hitsound volume rows --------------------------
$ws = $wb.Worksheets.Item("Miscellaneous")

$ws.Cells.Item(13, 1).Value = "SET_HITSOUND_VOLUME"
$ws.Cells.Item(13, 2).Value = "Set the overall hitsound volume"
$ws.Cells.Item(13, 3).Value = "전체적인 힛사운드 음량 조정하기"

$ws.Cells.Item(14, 1).Value = "CURRENT_HITSOUND_VOLUME"
$ws.Cells.Item(14, 2).Value = "Volume:"
$ws.Cells.Item(14, 3).Value = "음량:"
